# Apply updated figures to the "Crédito disponível - Centralização" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("M11").Value = 920047.22
$ws.Range("N11").Value = 472753.46
$ws.Range("O11").Value = 461023.84

# Row 12
$ws.Range("N12").Value = 144482.74
$ws.Range("O12").Value = 144482.74

# Row 13
$ws.Range("O13").Value = 8250

# Row 16
$ws.Range("K16").Value = 213987.19

# Row 17
$ws.Range("N17").Value = 3149.5

# Row 18
$ws.Range("K18").Value = 1135.95

$wb.Save()
